$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 1.09
$ws.Range("N2").Value = 1.34
$ws.Range("P2").Value = 1.34
$ws.Range("Q2").Value = 1.42

# Row 4 updates
$ws.Range("F4").Value = 1.72
$ws.Range("J4").Value = 3.65
$ws.Range("K4").Value = 4
$ws.Range("Q4").Value = 2.08
